$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column C (the "Förändrad" date column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# Update every data row (2..lastRow) in column C from 45172 -> 45175 (serial dates)
$rng = $ws.Range("C2:C$lastRow")
$rng.Value = 45175
